# Fruta / hortaliza, semanal
#
# The weekly refresh re-sorted the 8 data rows (rows 2-9) of the "Haba"
# sheet by date; every field in a data row travels together, so this is a
# row-level permutation (not independent cell tweaks). Target order,
# expressed as "new row <- old row":
#   2<-6  3<-7  4<-8  5<-9  6<-5  7<-3  8<-2  9<-4
#
# Only columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio
# maximo), M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg)
# actually differ row-to-row - A,B,C,E,F,G,H,I,N,Q,R are constant across
# all 8 rows, so nothing there needs to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 (= old row 6)
$ws.Range("D2").Value = 44159
$ws.Range("J2").Value = 42
$ws.Range("K2").Value = 6500
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 6738
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 270

# New row 3 (= old row 7)
$ws.Range("D3").Value = 44161
$ws.Range("J3").Value = 53
$ws.Range("K3").Value = 6500
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 6764
$ws.Range("O3").Value = "Región de O'Higgins"
$ws.Range("P3").Value = 271

# New row 4 (= old row 8)
$ws.Range("D4").Value = 44448
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 580

# New row 5 (= old row 9)
$ws.Range("D5").Value = 44167
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 8500
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 340

# New row 6 (= old row 5)
$ws.Range("D6").Value = 44165
$ws.Range("J6").Value = 38
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 8500
$ws.Range("M6").Value = 8263
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 331

# New row 7 (= old row 3)
$ws.Range("D7").Value = 44166
$ws.Range("J7").Value = 56
$ws.Range("K7").Value = 7500
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7804
$ws.Range("O7").Value = "Región de O'Higgins"
$ws.Range("P7").Value = 312

# New row 8 (= old row 2)
$ws.Range("D8").Value = 44160
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 6500
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 6688
$ws.Range("O8").Value = "Región de O'Higgins"
$ws.Range("P8").Value = 268

# New row 9 (= old row 4)
$ws.Range("D9").Value = 44162
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 7562
$ws.Range("O9").Value = "Región de O'Higgins"
$ws.Range("P9").Value = 302
